$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing source citation (currently at A23/A24) and move it
# down to A29/A30 to make room for the new MSME-definition table.
$sourceName = $ws.Range("A23").Value()
$sourceCite = $ws.Range("A24").Value()

$ws.Range("A23").Value = ""
$ws.Range("A24").Value = ""

# New table header (row 20), styled like the other header rows (11 / 16)
$ws.Range("B20").Value = "Number of employees"
$ws.Range("C20").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D20").Value = "Turnover (local currency, unless noted otherwise)"

# Data rows 21-24
$ws.Range("A21").Value = "Micro"
$ws.Range("B21").Value = "1-20"
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""

$ws.Range("A22").Value = "Small"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = "< limits set at the Global Unique Contribution (CGU) by the General Tax Code"

$ws.Range("A23").Value = "Medium"
$ws.Range("B23").Value = "21-250"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = "<= F CFA 5 Billionlion"

$ws.Range("A24").Value = "Large"
$ws.Range("B24").Value = ">250"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = "> F CFA 5 Billionlion"

# Restore the source citation further down the sheet
$ws.Range("A29").Value = $sourceName
$ws.Range("A30").Value = $sourceCite
